# Update the "valueObject" sheet's test-case field table (rows 27-31)
# to reflect the new Input/Expect value-object layout, and move the
# active selection to D32, per the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("valueObject")
$ws.Activate()

# --- Header block (class name / package / description) ---
$ws.Range("C6").Value = "BlancoRestAutotestTestCaseData"
$ws.Range("C7").Value = "blanco.restautotest.valueobject"
$ws.Range("C8").Value = "テストケースに記載されたデータを格納するためのクラスです。"

# --- Row 27: input field ---
$ws.Range("B27").Value = "input"
$ws.Range("C27").Value = "blanco.restgenerator.valueobject.ApiTelegram"
$ws.Range("D27").Value = ""
$ws.Range("E27").Value = "テストケースの入力値を格納します。"

# --- Row 28: inputId field ---
$ws.Range("B28").Value = "inputId"
$ws.Range("C28").Value = "java.lang.String"
$ws.Range("D28").Value = ""
$ws.Range("E28").Value = "電文クラスの正式名を格納します。"

# --- Row 29: expect field ---
$ws.Range("B29").Value = "expect"
$ws.Range("C29").Value = "blanco.restgenerator.valueobject.ApiTelegram"
$ws.Range("D29").Value = ""
$ws.Range("E29").Value = "テストケースの出力値を格納します。"

# --- Row 30: expectId field ---
$ws.Range("B30").Value = "expectId"
$ws.Range("C30").Value = "java.lang.String"
$ws.Range("D30").Value = ""
$ws.Range("E30").Value = "電文クラスの正式名を格納します。"

# --- Row 31: new caseId field (was blank) ---
$ws.Range("A31").Formula = "=A30+1"
$ws.Range("B31").Value = "caseId"
$ws.Range("C31").Value = "java.lang.String"
$ws.Range("D31").Value = ""
$ws.Range("E31").Value = "テストケースIDを格納します。"

# --- Move the active selection like the author left it ---
$ws.Range("D32").Select()
